$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 291, shifting existing rows 291.. down to 292..
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new data record
$ws.Cells.Item(291, 1).Value = 3
$ws.Cells.Item(291, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(291, 3).Value = "Coquimbo"
$ws.Cells.Item(291, 4).Value = 45215
$ws.Cells.Item(291, 5).Value = 5
$ws.Cells.Item(291, 6).Value = "Fruta"
$ws.Cells.Item(291, 7).Value = 100101
$ws.Cells.Item(291, 8).Value = "Berries"
$ws.Cells.Item(291, 9).Value = 100101001
$ws.Cells.Item(291, 10).Value = "Arándano (blue)"
$ws.Cells.Item(291, 11).Value = "Sin especificar"
$ws.Cells.Item(291, 12).Value = "Primera"
$ws.Cells.Item(291, 13).Value = 56
$ws.Cells.Item(291, 14).Value = 13000
$ws.Cells.Item(291, 15).Value = 13000
$ws.Cells.Item(291, 16).Value = 13000
$ws.Cells.Item(291, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(291, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(291, 19).Value = 6500
$ws.Cells.Item(291, 20).Value = 2

# Ensure the date cell keeps the same style/number format as the other date cells in column D
$ws.Cells.Item(291, 4).NumberFormat = $ws.Cells.Item(292, 4).NumberFormat
